$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Student DB register-dialog strings as rows 13-17 in columns B:C
$data = @(
    @("STR_NEW", "New"),
    @("STR_STUDENT_NAME_LIST", "Student List"),
    @("STR_NAME", "Name"),
    @("STR_BIRTH_PLACE", "Place of birth"),
    @("STR_BIRTH_DATE", "Date of birth")
)

$row = 13
foreach ($pair in $data) {
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $row = $row + 1
}

# Update the selection to match the next empty row, as in the diff
$ws.Range("C18").Select()
